$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New account-statement data: rows 16-54 of the worker/period table.
# Each worker's rows are now grouped together (periods 2312 -> 2306, newest
# first), replacing the old layout where rows were grouped by period with
# each worker appearing once per period block. Column B (doc type "CC") and
# column G (1160000) are unchanged for every row.
$rows = @(
    @{Row=16; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2312"; F=25333},
    @{Row=17; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2311"; F=46400},
    @{Row=18; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2310"; F=46400},
    @{Row=19; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2309"; F=46400},
    @{Row=20; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2308"; F=46400},
    @{Row=21; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2307"; F=46400},
    @{Row=22; C="1052973491"; D="RESLY RODRIGUEZ PALENCIA"; E="2306"; F=46400},
    @{Row=23; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2312"; F=25333},
    @{Row=24; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2311"; F=46400},
    @{Row=25; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2310"; F=46400},
    @{Row=26; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2309"; F=46400},
    @{Row=27; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2308"; F=46400},
    @{Row=28; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2307"; F=46400},
    @{Row=29; C="1143339170"; D="EDWIN CABARCAS ARIZA"; E="2306"; F=46400},
    @{Row=30; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2312"; F=28728},
    @{Row=31; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2311"; F=46400},
    @{Row=32; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2310"; F=46400},
    @{Row=33; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2309"; F=46400},
    @{Row=34; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2308"; F=46400},
    @{Row=35; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2307"; F=46400},
    @{Row=36; C="92071278"; D="NAIN URIEL MENDOZA MEDRANO"; E="2306"; F=46400},
    @{Row=37; C="1193456123"; D="ERIKA PATRICIA CASTRO ZABALETA"; E="2309"; F=46400},
    @{Row=38; C="1193456123"; D="ERIKA PATRICIA CASTRO ZABALETA"; E="2308"; F=46400},
    @{Row=39; C="1193456123"; D="ERIKA PATRICIA CASTRO ZABALETA"; E="2307"; F=46400},
    @{Row=40; C="1193456123"; D="ERIKA PATRICIA CASTRO ZABALETA"; E="2306"; F=46400},
    @{Row=41; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2312"; F=29387},
    @{Row=42; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2311"; F=46400},
    @{Row=43; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2310"; F=46400},
    @{Row=44; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2309"; F=46400},
    @{Row=45; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2308"; F=46400},
    @{Row=46; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2307"; F=46400},
    @{Row=47; C="1007275438"; D="MILEIDIS POLANCO PADILLA"; E="2306"; F=46400},
    @{Row=48; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2312"; F=29387},
    @{Row=49; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2311"; F=46400},
    @{Row=50; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2310"; F=46400},
    @{Row=51; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2309"; F=46400},
    @{Row=52; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2308"; F=46400},
    @{Row=53; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2307"; F=46400},
    @{Row=54; C="1052954248"; D="LORENA PATRICIA ROMERO CASTRO"; E="2306"; F=46400}
)

foreach ($row in $rows) {
    $r = $row.Row
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F
    $ws.Range("G$r").Value = 1160000
}
